$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value2 = $value
    $r.Style = "Normal"
}

Set-TextValue "D2" '37.723.65'
Set-TextValue "E2" '  -1.55%  '
Set-TextValue "D3" '2.077.06'
Set-TextValue "E3" '  -2.23%  '
Set-TextValue "E4" '  -0.11%  '
Set-TextValue "D5" '233.60'
Set-TextValue "E5" '  -1.00%  '
Set-TextValue "D6" '0.623'
Set-TextValue "E6" '  -0.76%  '
Set-TextValue "D7" '58.74'
Set-TextValue "E7" '  +0.45%  '
Set-TextValue "E8" '  -0.11%  '
Set-TextValue "E9" '  +0.28%  '
Set-TextValue "D10" '0.0784'
Set-TextValue "E10" '  -0.07%  '
Set-TextValue "E11" '  +3.29%  '
Set-TextValue "D12" '2.381.84'
Set-TextValue "E12" '  -2.30%  '
Set-TextValue "D13" '14.82'
Set-TextValue "E13" '  +1.27%  '
Set-TextValue "D14" '20.97'
Set-TextValue "E14" '  -3.25%  '
Set-TextValue "D15" '0.774'
Set-TextValue "E15" '  -2.66%  '
Set-TextValue "E16" '  +1.90%  '
Set-TextValue "D17" '2.065.51'
Set-TextValue "E17" '  -2.71%  '
Set-TextValue "D18" '37.651.78'
Set-TextValue "E18" '  -1.44%  '
Set-TextValue "D19" '6.12'
Set-TextValue "E19" '  -1.82%  '
Set-TextValue "D20" '71.42'
Set-TextValue "E20" '  +0.97%  '
Set-TextValue "D21" '0.0₃0835'
Set-TextValue "E21" '  +0.69%  '
Set-TextValue "D22" '228.48'
Set-TextValue "E22" '  -0.61%  '
Set-TextValue "E23" '  -0.08%  '
Set-TextValue "D24" '2.41'
Set-TextValue "E24" '  -0.42%  '
Set-TextValue "E25" '  -2.70%  '
Set-TextValue "D26" '171.40'
Set-TextValue "E26" '  +1.28%  '
Set-TextValue "D27" '9.08'
Set-TextValue "E28" '  -1.79%  '
Set-TextValue "D29" '19.51'
Set-TextValue "E29" '  -0.74%  '
Set-TextValue "D30" '1.40'
Set-TextValue "E30" '  -2.54%  '
Set-TextValue "E31" '  +1.70%  '
Set-TextValue "E32" '  +0.35%  '
Set-TextValue "E33" '  +0.65%  '
Set-TextValue "D34" '4.66'
Set-TextValue "E34" '  +1.15%  '
Set-TextValue "E35" '  -5.61%  '
Set-TextValue "D36" '1.83'
Set-TextValue "E36" '  -0.40%  '
Set-TextValue "E37" '  -2.98%  '
Set-TextValue "E38" '  +0.05%  '
Set-TextValue "E39" '  -2.69%  '
Set-TextValue "D40" '0.0974'
Set-TextValue "E40" '  -2.99%  '
Set-TextValue "D41" '99.65'
Set-TextValue "E41" '  +2.09%  '
Set-TextValue "E42" '  -2.44%  '
Set-TextValue "E43" '  -0.02%  '
Set-TextValue "D44" '16.64'
Set-TextValue "E44" '  +5.20%  '
Set-TextValue "D45" '1.436.94'
Set-TextValue "E45" '  -1.83%  '
Set-TextValue "E46" '  -1.27%  '
Set-TextValue "D47" '4.21'
Set-TextValue "E47" '  +2.51%  '
Set-TextValue "E48" '  -0.99%  '
Set-TextValue "E49" '  +0.74%  '
Set-TextValue "E50" '  -1.67%  '
Set-TextValue "D51" '2.266.75'
Set-TextValue "E51" '  -2.36%  '
